$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.329339861869812
$ws.Range("B1").Value = 1.545700430870056
$ws.Range("C1").Value = 1.976040959358215
$ws.Range("D1").Value = 1.931671142578125
$ws.Range("E1").Value = 1.607928276062012
